# Updated cryptos list values (price + volume columns) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values are plain decimals (e.g. "593.07") that Excel would
# otherwise auto-convert to numbers on assignment. The source data stores
# these as literal text, so format as Text before writing, then restore the
# default "Normal" style afterwards so no visible/style change is left behind.
$textCells = @("D4", "D5", "D6", "D11", "D14", "D21", "D22", "D25", "D28", "D29", "D30", "D32", "D34", "D40", "D41", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '65.084.95'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '3.526.19'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('D5').Value = '593.07'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').Value = '134.24'
$ws.Range('E6').Value = '  -0.61%  '
$ws.Range('D7').Value = '3.525.47'
$ws.Range('E7').Value = '  -0.60%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -0.34%  '
$ws.Range('E10').Value = '  +2.11%  '
$ws.Range('D11').Value = '7.13'
$ws.Range('E11').Value = '  +3.45%  '
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').Value = '4.121.21'
$ws.Range('E13').Value = '  -0.71%  '
$ws.Range('D14').Value = '27.72'
$ws.Range('E14').Value = '  +2.66%  '
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('E16').Value = '  +0.59%  '
$ws.Range('D17').Value = '3.525.79'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('D18').Value = '65.046.84'
$ws.Range('E18').Value = '  +0.61%  '
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').Value = '5.69'
$ws.Range('E21').Value = '  -2.07%  '
$ws.Range('D22').Value = '392.55'
$ws.Range('E22').Value = '  +1.73%  '
$ws.Range('E23').Value = '  +0.93%  '
$ws.Range('B24').Value = 'WrappedeETH'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D24').Value = '3.666.90'
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '74.83'
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -2.71%  '
$ws.Range('B28').Value = 'Fetch.AI'
$ws.Range('C28').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D28').Value = '1.61'
$ws.Range('E28').Value = '  +10.96%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '7.75'
$ws.Range('E29').Value = '  +1.79%  '
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('D32').Value = '8.36'
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('D33').Value = '3.532.19'
$ws.Range('E33').Value = '  -0.60%  '
$ws.Range('D34').Value = '24.16'
$ws.Range('E34').Value = '  +0.88%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  +0.42%  '
$ws.Range('E38').Value = '  +1.00%  '
$ws.Range('E39').Value = '  +2.85%  '
$ws.Range('D40').Value = '168.38'
$ws.Range('E40').Value = '  -0.82%  '
$ws.Range('D41').Value = '0.0811'
$ws.Range('E41').Value = '  +1.02%  '
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('E43').Value = '  +5.97%  '
$ws.Range('D44').Value = '25.79'
$ws.Range('E44').Value = '  -4.33%  '
$ws.Range('D45').Value = '42.98'
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = '4.45'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('D48').Value = '1.67'
$ws.Range('E48').Value = '  +1.82%  '
$ws.Range('D49').Value = '6.91'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').Value = '2.430.86'
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('D51').Value = '0.912'
$ws.Range('E51').Value = '  +6.76%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
